# "Made delete work and started on sorts"
#
# Rebuilds Sheet1 as an 8-column (A:H) x 6-row (1:6) table:
#   - col A is a fresh 0-based row index (bold/bordered header style, no header text)
#   - cols B/C are new "Unnamed: 0.1" / "Unnamed: 0" index columns carried over
#     from a prior dataframe export, headered in the same bold style as D1:H1
#   - cols D:H are the original numeric-named columns (0,1,2,3,4), now holding
#     five re-sorted "Transaction N" rows instead of the original two test rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# $null marks a cell that must stay empty (A1 - blank corner above the index column).
$data = @(
  @($null, "Unnamed: 0.1", "Unnamed: 0", 0, 1, 2, 3, 4),
  @(0, 0, 1, 5, "Transaction 1", 100, 10, 1),
  @(1, 4, 5, 1, "Transaction 5", 500, 50, 5),
  @(2, 3, 4, 2, "Transaction 4", 400, 40, 4),
  @(3, 2, 3, 3, "Transaction 3", 300, 30, 3),
  @(4, 1, 2, 4, "Transaction 2", 200, 20, 2)
)

for ($ri = 0; $ri -lt $data.Length; $ri++) {
  $row = $data[$ri]
  $r = $ri + 1
  for ($ci = 0; $ci -lt $row.Length; $ci++) {
    $c = $ci + 1
    $val = $row[$ci]
    if ($null -ne $val) {
      $ws.Cells.Item($r, $c).Value = $val
    }
  }
}

# The sheet already carries the bold/bordered "header" style on B1:F1 and A2:A3
# from before this edit. Stamp that same style onto the newly-added header
# cells (G1:H1) and the newly-added index rows (A4:A6) so the whole extended
# header row and index column are consistently styled.
$ws.Range("D1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$ws.Range("A4:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
